{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Goal: after the \"Code structuring\" heading paragraph, insert a new body\n// paragraph discussing Socket.IO's asynchronous / event based nature\n// (with a red \"TODO: crossref\" placeholder), followed by one empty\n// paragraph - matching the target diff.\n\n// 1. Locate the \"Code structuring\" heading paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet headingPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Code structuring\") {\n    headingPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!headingPara) {\n  throw new Error('Could not find the \"Code structuring\" paragraph.');\n}\n\n// 2. Build the new paragraphs as raw WordprocessingML, wrapped in the\n// flat-OPC package envelope insertOoxml() expects. Using OOXML (rather than\n// insertParagraph/insertText) lets us emit the exact run layout - including\n// the color-only run carrying the red \"TODO: crossref\" placeholder, its\n// proofErr spell-check markers, and a following paragraph that does NOT\n// inherit the Heading 3 paragraph style (the engine's insertParagraph()\n// APIs always copy the anchor paragraph's pPr, which would incorrectly\n// leave the new body text styled as Heading 3).\nconst newParagraphsXml =\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">Perhaps the most unfamiliar aspect of JavaScript compared to other languages, is the fact that it is asynchronous. A common pitfall for JavaScript frameworks is to only provide developers with synchronous tools to use. With Socket.IO this is, thankfully, not the case. Socket.IO follows the WebSockets protocol tightly as it provides and event based architecture. While the WebSockets API only provides </w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>few, standard events (</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">TODO: </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>crossref</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>)</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>, Socket.IO lets you used self named events in addition to the standard WebSocket API events.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr></w:pPr></w:p>';\n\nconst flatOpcXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newParagraphsXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// 3. Insert the new paragraphs right after the \"Code structuring\" heading.\nconst endOfHeading = headingPara.getRange(\"End\");\nendOfHeading.insertOoxml(flatOpcXml, \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Goal: after the \"Code structuring\" heading paragraph, insert a new body\n# paragraph discussing Socket.IO's asynchronous / event based nature\n# (with a red \"TODO: crossref\" placeholder), followed by one empty\n# paragraph - matching the target diff.\n\n$d = $word.ActiveDocument\n\n# 1. Locate the \"Code structuring\" heading paragraph and remember its\n# 1-based index in $d.Paragraphs so we can re-fetch fresh Range objects\n# after each mutation (ranges captured before an insert do not \"move\" in\n# this host, so re-resolving by index keeps us anchored correctly).\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n  $i = $i + 1\n  if ($p.Range.Text.Trim() -eq \"Code structuring\") {\n    $targetIndex = $i\n    break\n  }\n}\n\nif ($targetIndex -lt 0) {\n  throw 'Could not find the \"Code structuring\" paragraph.'\n}\n\n# 2. Create a brand new, empty paragraph right after it.\n$headingRange = $d.Paragraphs.Item($targetIndex).Range\n$headingRange.Collapse(0)\n$headingRange.InsertParagraphAfter()\n\n# 3. Re-fetch the freshly created (now separate) paragraph and fill it in\n# via raw WordprocessingML, wrapped in the flat-OPC package envelope\n# InsertXML() expects. Using OOXML (rather than typing plain text) lets us\n# emit the exact run layout - including the color-only run carrying the\n# red \"TODO: crossref\" placeholder and its proofErr spell-check markers -\n# and a following empty paragraph, without inheriting the Heading 3 style\n# that a plain InsertParagraphAfter() copies from its anchor paragraph.\n$newParagraphsXml = (\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">Perhaps the most unfamiliar aspect of JavaScript compared to other languages, is the fact that it is asynchronous. A common pitfall for JavaScript frameworks is to only provide developers with synchronous tools to use. With Socket.IO this is, thankfully, not the case. Socket.IO follows the WebSockets protocol tightly as it provides and event based architecture. While the WebSockets API only provides </w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>few, standard events (</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">TODO: </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>crossref</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>)</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr><w:t>, Socket.IO lets you used self named events in addition to the standard WebSocket API events.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\" w:eastAsia=\"en-US\"/></w:rPr></w:pPr></w:p>'\n)\n\n$flatOpcXml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>$newParagraphsXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n\"@\n\n$newParaRange = $d.Paragraphs.Item($targetIndex + 1).Range\n$newParaRange.InsertXML($flatOpcXml)\n"}
